$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "publish date" labels in row 9 (quarterly report revision tags) ---
$ws.Range("I9").Value = "1402-01-28 (5)"
$ws.Range("J9").Value = "1402-01-28 (8)"
$ws.Range("M9").Value = "1402-01-28 (3)"

# --- Updated figures for the most-recently revised quarter (column J) per the
#     new read_price algorithm ---
$ws.Range("J11").Value = 91319    # فروش (Sales)
$ws.Range("J12").Value = -76236   # بهای تمام شده کالای فروش رفته (COGS)
$ws.Range("J13").Value = 15083    # سود (زیان) ناخالص (Gross profit)
$ws.Range("J17").Value = 12842    # سود (زیان) عملیاتی (Operating profit)
$ws.Range("J20").Value = 12617    # سود خالص عملیات قبل از مالیات
$ws.Range("J22").Value = 10863    # سود (زیان) خالص عملیات در حال تداوم
$ws.Range("J24").Value = 10863    # سود (زیان) خالص

# --- EPS row: previously "-" placeholders, now populated with 0 ---
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("M25").Value = 0

# --- سرمایه (Capital) row: previously "-" placeholders, now populated with figures ---
$ws.Range("I26").Value = 7025
$ws.Range("J26").Value = 7261
$ws.Range("M26").Value = 5757
